$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update column B (Taxonsorteringsordning) for rows 2-10: 57880 -> 57884
$ws.Range("B2:B10").Value = 57884

# Rows 8 and 9 have their data swapped (except column B which is already handled above).
# Capture current ("before") values for the columns that differ between the two rows.
$A8  = $ws.Range("A8").Value()
$M8  = $ws.Range("M8").Value()
$Q8  = $ws.Range("Q8").Value()
$R8  = $ws.Range("R8").Value()
$Z8  = $ws.Range("Z8").Value()
$AB8 = $ws.Range("AB8").Value()
$AC8 = $ws.Range("AC8").Value()

$A9  = $ws.Range("A9").Value()
$M9  = $ws.Range("M9").Value()
$Q9  = $ws.Range("Q9").Value()
$R9  = $ws.Range("R9").Value()
$Z9  = $ws.Range("Z9").Value()
$AB9 = $ws.Range("AB9").Value()
$AC9 = $ws.Range("AC9").Value()

# Write row 9's former values into row 8
$ws.Range("A8").Value = $A9
$ws.Range("M8").Value = $M9
$ws.Range("Q8").Value = $Q9
$ws.Range("R8").Value = $R9
$ws.Range("Z8").Value = $Z9
$ws.Range("AB8").Value = $AB9
$ws.Range("AC8").Value = $AC9

# Write row 8's former values into row 9
$ws.Range("A9").Value = $A8
$ws.Range("M9").Value = $M8
$ws.Range("Q9").Value = $Q8
$ws.Range("R9").Value = $R8
$ws.Range("Z9").Value = $Z8
$ws.Range("AB9").Value = $AB8
$ws.Range("AC9").Value = $AC8
